$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The update swaps the "Primera" quality pair (rows 2-3) and the
# "Segunda" quality pair (rows 4-5) for columns D, J, K, L, M, N, P, Q
# while leaving the identifying columns (A, B, C, E, F, G, H, I, O, R)
# untouched.

# Row 2 <- (former) row 3 values
$ws.Range("D2").Value = 44623
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1900
$ws.Range("N2").Value = '$/paquete'
$ws.Range("P2").Value = 1900
$ws.Range("Q2").Value = 1

# Row 3 <- (former) row 2 values
$ws.Range("D3").Value = 44267
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = 1650
$ws.Range("N3").Value = '$/docena de matas'
$ws.Range("P3").Value = 275
$ws.Range("Q3").Value = 6

# Row 4 <- (former) row 5 values
$ws.Range("D4").Value = 44370
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1200
$ws.Range("M4").Value = 1080
$ws.Range("N4").Value = '$/docena de matas'
$ws.Range("P4").Value = 180
$ws.Range("Q4").Value = 6

# Row 5 <- (former) row 4 values
$ws.Range("D5").Value = 44377
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2800
$ws.Range("M5").Value = 2364
$ws.Range("N5").Value = '$/docena de matas'
$ws.Range("P5").Value = 394
$ws.Range("Q5").Value = 6
